$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 2.1
$ws.Range("L2").Value = 4.6
$ws.Range("T2").Value = 2.2
$ws.Range("U2").Value = 1.78
$ws.Range("W2").Value = 6.4
$ws.Range("X2").Value = 9.75
$ws.Range("AA2").Value = 18
$ws.Range("AH2").Value = 9.75
$ws.Range("AK2").Value = 75
$ws.Range("AL2").Value = 45
$ws.Range("AM2").Value = 40

# Row 3
$ws.Range("G3").Value = 1.82
$ws.Range("I3").Value = 4.1
$ws.Range("J3").Value = 2.35
$ws.Range("K3").Value = 2.12
$ws.Range("L3").Value = 4.5
$ws.Range("U3").Value = 1.85
$ws.Range("V3").Value = 1.75
$ws.Range("W3").Value = 6.4
$ws.Range("X3").Value = 8.25
$ws.Range("Y3").Value = 8.25
$ws.Range("Z3").Value = 15
$ws.Range("AA3").Value = 15.5
$ws.Range("AE3").Value = 16.5
$ws.Range("AF3").Value = 90
$ws.Range("AH3").Value = 10.25
$ws.Range("AI3").Value = 22
$ws.Range("AJ3").Value = 14
$ws.Range("AK3").Value = 65
$ws.Range("AL3").Value = 45
$ws.Range("AM3").Value = 50

# Row 4
$ws.Range("H4").Value = 5.6
$ws.Range("I4").Value = 11.5
$ws.Range("K4").Value = 2.6
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 2.25
$ws.Range("R4").Value = 1.5
$ws.Range("U4").Value = 2.12
$ws.Range("V4").Value = 1.57
$ws.Range("X4").Value = 5.9
$ws.Range("Z4").Value = 6.8
$ws.Range("AB4").Value = 32
$ws.Range("AI4").Value = 90
$ws.Range("AJ4").Value = 37
$ws.Range("AK4").Value = 400
$ws.Range("AM4").Value = 120

# Row 5
$ws.Range("G5").Value = 2.42
$ws.Range("H5").Value = 2.9
$ws.Range("J5").Value = 3.1
$ws.Range("L5").Value = 3.7
$ws.Range("M5").Value = 1.5
$ws.Range("N5").Value = 2.27
$ws.Range("O5").Value = 2.4
$ws.Range("P5").Value = 1.44
$ws.Range("Q5").Value = 4.15
$ws.Range("R5").Value = 1.15
$ws.Range("S5").Value = 1.53
$ws.Range("U5").Value = 2.02
$ws.Range("V5").Value = 1.62
$ws.Range("W5").Value = 6.1
$ws.Range("X5").Value = 10.5
$ws.Range("Y5").Value = 10
$ws.Range("AA5").Value = 25
$ws.Range("AB5").Value = 45
$ws.Range("AC5").Value = 6.4
$ws.Range("AE5").Value = 18
$ws.Range("AM5").Value = 55

# Row 6
$ws.Range("H6").Value = 3.25
$ws.Range("I6").Value = 2.67
$ws.Range("J6").Value = 3.05
$ws.Range("K6").Value = 2.07
$ws.Range("L6").Value = 3.25
$ws.Range("M6").Value = 1.36
$ws.Range("N6").Value = 2.67
$ws.Range("O6").Value = 2.05
$ws.Range("P6").Value = 1.6
$ws.Range("Q6").Value = 3.4
$ws.Range("R6").Value = 1.22
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.5
$ws.Range("U6").Value = 1.83
$ws.Range("V6").Value = 1.78
$ws.Range("W6").Value = 7.3
$ws.Range("X6").Value = 11.25
$ws.Range("AA6").Value = 22
$ws.Range("AB6").Value = 35
$ws.Range("AC6").Value = 8.5
$ws.Range("AD6").Value = 6.3
$ws.Range("AE6").Value = 16
$ws.Range("AF6").Value = 90
$ws.Range("AG6").Value = 800
$ws.Range("AH6").Value = 7.6
$ws.Range("AJ6").Value = 10.25
$ws.Range("AL6").Value = 24
$ws.Range("AM6").Value = 37

# Row 7
$ws.Range("H7").Value = 3.65
$ws.Range("M7").Value = 1.27
$ws.Range("N7").Value = 3.1
$ws.Range("O7").Value = 1.8
$ws.Range("P7").Value = 1.82
$ws.Range("Q7").Value = 2.85
$ws.Range("R7").Value = 1.32
$ws.Range("U7").Value = 1.75
$ws.Range("V7").Value = 1.85
$ws.Range("W7").Value = 7.2
$ws.Range("X7").Value = 8.5
$ws.Range("Z7").Value = 14.5
$ws.Range("AB7").Value = 26
$ws.Range("AD7").Value = 7.1
$ws.Range("AE7").Value = 15.5
$ws.Range("AF7").Value = 70
$ws.Range("AG7").Value = 600

# Row 10
$ws.Range("G10").Value = 2.38
$ws.Range("I10").Value = 3
$ws.Range("M10").Value = 1.3
$ws.Range("O10").Value = 2
$ws.Range("P10").Value = 1.8
$ws.Range("Q10").Value = 3.4
$ws.Range("R10").Value = 1.3
$ws.Range("S10").Value = 1.4
$ws.Range("T10").Value = 2.75
$ws.Range("U10").Value = 1.73
$ws.Range("V10").Value = 2
$ws.Range("X10").Value = 12
$ws.Range("AC10").Value = 9.5
$ws.Range("AF10").Value = 41
$ws.Range("AG10").Value = 201
$ws.Range("AL10").Value = 23
$ws.Range("AN10").Value = 1.06

# Row 11
$ws.Range("G11").Value = 2.15
$ws.Range("H11").Value = 2.9
$ws.Range("I11").Value = 3.5
$ws.Range("J11").Value = 2.88
$ws.Range("K11").Value = 1.91
$ws.Range("L11").Value = 4.33
$ws.Range("M11").Value = 1.5
$ws.Range("O11").Value = 2.5
$ws.Range("P11").Value = 1.5
$ws.Range("R11").Value = 1.17
$ws.Range("S11").Value = 1.57
$ws.Range("T11").Value = 2.25
$ws.Range("U11").Value = 2.2
$ws.Range("V11").Value = 1.62
$ws.Range("X11").Value = 9
$ws.Range("Y11").Value = 10
$ws.Range("Z11").Value = 19
$ws.Range("AA11").Value = 21
$ws.Range("AC11").Value = 6.5
$ws.Range("AE11").Value = 19
$ws.Range("AF11").Value = 67
$ws.Range("AH11").Value = 8.5
$ws.Range("AI11").Value = 17
$ws.Range("AM11").Value = 41
$ws.Range("AN11").Value = 1.1
$ws.Range("AO11").Value = 6.5
$ws.Range("AP11").Value = 1.93
$ws.Range("AQ11").Value = 1.93

# Row 12
$ws.Range("G12").Value = 3
$ws.Range("H12").Value = 2.85
$ws.Range("I12").Value = 2.45
$ws.Range("J12").Value = 3.75
$ws.Range("L12").Value = 3.25
$ws.Range("M12").Value = 1.44
$ws.Range("N12").Value = 2.63
$ws.Range("O12").Value = 2.4
$ws.Range("P12").Value = 1.53
$ws.Range("Q12").Value = 4.5
$ws.Range("R12").Value = 1.18
$ws.Range("S12").Value = 1.53
$ws.Range("T12").Value = 2.38
$ws.Range("U12").Value = 2
$ws.Range("V12").Value = 1.73
$ws.Range("W12").Value = 7.5
$ws.Range("X12").Value = 13
$ws.Range("Y12").Value = 12
$ws.Range("AA12").Value = 29
$ws.Range("AI12").Value = 11
$ws.Range("AJ12").Value = 11
$ws.Range("AK12").Value = 23
$ws.Range("AL12").Value = 23
$ws.Range("AN12").Value = 1.1
$ws.Range("AP12").Value = 1.83
$ws.Range("AQ12").Value = 2.03

# Row 13
$ws.Range("G13").Value = 2.77
$ws.Range("H13").Value = 3.15
$ws.Range("I13").Value = 2.37
$ws.Range("J13").Value = 3.35
$ws.Range("L13").Value = 3.05
$ws.Range("T13").Value = 2.7
$ws.Range("W13").Value = 10
$ws.Range("Y13").Value = 10
$ws.Range("AC13").Value = 7.4
$ws.Range("AE13").Value = 12.5
$ws.Range("AI13").Value = 12
$ws.Range("AK13").Value = 26
$ws.Range("AO13").Value = 7.4
